# Regenerate orders with updated distance/sizes.
# Applies the following global text substitutions to every string-valued
# cell in the worksheet's used range:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# (applied to whole cell values, which also fixes up composite values like
#  "Face08_D64_S25", "Face01_D64_S30_l.png", etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()

        if ($v -is [string]) {
            $newV = $v.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($newV -ne $v) {
                $cell.Value = $newV
            }
        }
    }
}
